$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix part rotations (orientation) for production:
#  - J1 (row 8): 90 -> -90
#  - J2 (row 9): -90 -> 90
#  - Q1 (row 10): -90 -> 90
$ws.Range("E8").Value = -90
$ws.Range("E9").Value = 90
$ws.Range("E10").Value = 90

# Re-apply the data rows' cell style so it matches the (already-defined,
# reused) style used elsewhere in the sheet: Arial 11, centered horizontally
# and vertically, no border/fill. Copying format from an existing cell that
# already carries this exact combination avoids creating brand new font /
# cell-format entries in the workbook's style table.
$fmtSource = $ws.Range("G1")
$dataRange = $ws.Range("A2:E15")
$fmtSource.Copy()
$dataRange.PasteSpecial(-4122)
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108

# Update the saved selection/active cell.
$ws.Range("E22").Select()
